# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Atomos_Profits workbook
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit* per leve row)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 378.57144
$ws.Range("J28").Value = 520
$ws.Range("L28").Value = 520
$ws.Range("N28").Value = -1490
$ws.Range("H98").Value = 2175.818
$ws.Range("I98").Value = 1315.6666
$ws.Range("J98").Value = 4019
$ws.Range("K98").Value = 1315.6666
$ws.Range("L98").Value = 4019
$ws.Range("M98").Value = 182.3334
$ws.Range("N98").Value = -7015
$ws.Range("H111").Value = 1850
$ws.Range("I111").Value = 1850
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 5550
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -2483
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 2175.818
$ws.Range("I122").Value = 1315.6666
$ws.Range("J122").Value = 4019
$ws.Range("K122").Value = 3946.9998
$ws.Range("L122").Value = 12057
$ws.Range("M122").Value = -1496.9998
$ws.Range("N122").Value = -16957
$ws.Range("H129").Value = 1036.42
$ws.Range("J129").Value = 1103
$ws.Range("L129").Value = 3309
$ws.Range("N129").Value = -13309
$ws.Range("H138").Value = 2538.3635
$ws.Range("J138").Value = 4052.577
$ws.Range("L138").Value = 12157.731
$ws.Range("N138").Value = -22437.731

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1112267.6
$ws.Range("I86").Value = 1726178.2
$ws.Range("J86").Value = 89083.336
$ws.Range("K86").Value = 1726178.2
$ws.Range("L86").Value = 89083.336
$ws.Range("M86").Value = -1725055.2
$ws.Range("N86").Value = -91329.336
$ws.Range("H89").Value = 1112267.6
$ws.Range("I89").Value = 1726178.2
$ws.Range("J89").Value = 89083.336
$ws.Range("K89").Value = 8630891
$ws.Range("L89").Value = 445416.68
$ws.Range("M89").Value = -8625275
$ws.Range("N89").Value = -456648.68
$ws.Range("H141").Value = 34247
$ws.Range("J141").Value = 27945.8
$ws.Range("L141").Value = 27945.8
$ws.Range("N141").Value = -38305.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2002711.5
$ws.Range("I31").Value = 3126629
$ws.Range("J31").Value = 4636.1113
$ws.Range("K31").Value = 3126629
$ws.Range("L31").Value = 4636.1113
$ws.Range("M31").Value = -3126334
$ws.Range("N31").Value = -5226.1113
$ws.Range("H34").Value = 2002711.5
$ws.Range("I34").Value = 3126629
$ws.Range("J34").Value = 4636.1113
$ws.Range("K34").Value = 3126629
$ws.Range("L34").Value = 4636.1113
$ws.Range("M34").Value = -3126427
$ws.Range("N34").Value = -5040.1113
$ws.Range("H41").Value = 3240.3333
$ws.Range("I41").Value = 1149.8334
$ws.Range("J41").Value = 7421.3335
$ws.Range("K41").Value = 1149.8334
$ws.Range("L41").Value = 7421.3335
$ws.Range("M41").Value = -721.8334
$ws.Range("N41").Value = -8277.333500000001
$ws.Range("H99").Value = 2863.5625
$ws.Range("I99").Value = 2027.3334
$ws.Range("K99").Value = 2027.3334
$ws.Range("M99").Value = -529.3334
$ws.Range("H105").Value = 2545
$ws.Range("I105").Value = 2431
$ws.Range("J105").Value = 2687.5
$ws.Range("K105").Value = 2431
$ws.Range("L105").Value = 2687.5
$ws.Range("M105").Value = -684
$ws.Range("N105").Value = -6181.5
$ws.Range("H107").Value = 1902.8182
$ws.Range("I107").Value = 604.4286
$ws.Range("J107").Value = 4175
$ws.Range("K107").Value = 604.4286
$ws.Range("L107").Value = 4175
$ws.Range("M107").Value = 1315.5714
$ws.Range("N107").Value = -8015
$ws.Range("H126").Value = 2863.5625
$ws.Range("I126").Value = 2027.3334
$ws.Range("K126").Value = 6082.0002
$ws.Range("M126").Value = -3612.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H101").Value = 2962.5
$ws.Range("J101").Value = 2962.5
$ws.Range("L101").Value = 8887.5
$ws.Range("N101").Value = -13755.5
$ws.Range("H102").Value = 2100
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 9000
$ws.Range("N102").Value = -13868
$ws.Range("H114").Value = 1172.7778
$ws.Range("I114").Value = 1242.6666
$ws.Range("J114").Value = 1158.8
$ws.Range("K114").Value = 3727.9998
$ws.Range("L114").Value = 3476.4
$ws.Range("M114").Value = -473.9998000000001
$ws.Range("N114").Value = -9984.4
$ws.Range("H136").Value = 2309.889
$ws.Range("I136").Value = 1612.7142
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 4838.142599999999
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = 261.8574000000008
$ws.Range("N136").Value = -24450

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5007.3335
$ws.Range("I70").Value = 5234.222
$ws.Range("J70").Value = 4326.6665
$ws.Range("K70").Value = 5234.222
$ws.Range("L70").Value = 4326.6665
$ws.Range("M70").Value = -4964.222
$ws.Range("N70").Value = -4866.6665
$ws.Range("H73").Value = 5007.3335
$ws.Range("I73").Value = 5234.222
$ws.Range("J73").Value = 4326.6665
$ws.Range("K73").Value = 5234.222
$ws.Range("L73").Value = 4326.6665
$ws.Range("M73").Value = -4298.222
$ws.Range("N73").Value = -6198.6665
$ws.Range("H102").Value = 129023.875
$ws.Range("I102").Value = 3270.6667
$ws.Range("K102").Value = 3270.6667
$ws.Range("M102").Value = -1648.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 32666.666
$ws.Range("I4").Value = 30000
$ws.Range("J4").Value = 33200
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 33200
$ws.Range("M4").Value = -29887
$ws.Range("N4").Value = -33426
$ws.Range("H7").Value = 1764.8125
$ws.Range("I7").Value = 951.7273
$ws.Range("J7").Value = 3553.6
$ws.Range("K7").Value = 951.7273
$ws.Range("L7").Value = 3553.6
$ws.Range("M7").Value = -839.7273
$ws.Range("N7").Value = -3777.6
$ws.Range("H25").Value = 59672
$ws.Range("J25").Value = 59672
$ws.Range("L25").Value = 59672
$ws.Range("N25").Value = -60132
$ws.Range("H28").Value = 32666.666
$ws.Range("I28").Value = 30000
$ws.Range("J28").Value = 33200
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 33200
$ws.Range("M28").Value = -29768
$ws.Range("N28").Value = -33664
$ws.Range("H37").Value = 32666.666
$ws.Range("I37").Value = 30000
$ws.Range("J37").Value = 33200
$ws.Range("K37").Value = 30000
$ws.Range("L37").Value = 33200
$ws.Range("M37").Value = -29893
$ws.Range("N37").Value = -33414
$ws.Range("H126").Value = 1764.8125
$ws.Range("I126").Value = 951.7273
$ws.Range("J126").Value = 3553.6
$ws.Range("K126").Value = 2855.1819
$ws.Range("L126").Value = 10660.8
$ws.Range("M126").Value = -385.1819
$ws.Range("N126").Value = -15600.8
$ws.Range("H132").Value = 2400.7
$ws.Range("I132").Value = 1635.0555
$ws.Range("K132").Value = 4905.166499999999
$ws.Range("M132").Value = -2375.166499999999
$ws.Range("H135").Value = 29898.215
$ws.Range("J135").Value = 29898.215
$ws.Range("L135").Value = 29898.215
$ws.Range("N135").Value = -40038.215

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 22751.25
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 25858.572
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 25858.572
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -26194.572
$ws.Range("H122").Value = 502519
$ws.Range("I122").Value = 668645.4
$ws.Range("J122").Value = 4139.8
$ws.Range("K122").Value = 2005936.2
$ws.Range("L122").Value = 12419.4
$ws.Range("M122").Value = -2003486.2
$ws.Range("N122").Value = -17319.4
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
